# Data update: increment "Pagos" (F) and "Inscrições homologadas" (H) by 1
# for the rows below, keeping "Isenções deferidas" (G) unchanged
# (H = F + G in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$rows = @(15, 18, 41, 48, 62, 63, 83)

foreach ($r in $rows) {
    $fCell = $ws.Cells.Item($r, 6)
    $hCell = $ws.Cells.Item($r, 8)
    $fCell.Value = $fCell.Value() + 1
    $hCell.Value = $hCell.Value() + 1
}
